$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 43-56 (rows shift: old data moves down by 2,
# and two new price records are inserted as the new rows 43-44) ---
# Row 43
$ws.Range("D43").Value = 44460
$ws.Range("L43").Value = 'Especial'
$ws.Range("M43").Value = 20

# Row 44
$ws.Range("D44").Value = 44460
$ws.Range("M44").Value = 60
$ws.Range("N44").Value = 2800
$ws.Range("O44").Value = 3000
$ws.Range("P44").Value = 2900
$ws.Range("S44").Value = 2900

# Row 45
$ws.Range("D45").Value = 44414
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 55
$ws.Range("N45").Value = 3500
$ws.Range("O45").Value = 3500
$ws.Range("P45").Value = 3500
$ws.Range("Q45").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S45").Value = 3500
$ws.Range("T45").Value = 1

# Row 46
$ws.Range("D46").Value = 44165
$ws.Range("M46").Value = 50
$ws.Range("N46").Value = 2300
$ws.Range("O46").Value = 2300
$ws.Range("P46").Value = 2300
$ws.Range("S46").Value = 2300

# Row 47
$ws.Range("D47").Value = 44427
$ws.Range("L47").Value = 'Especial'
$ws.Range("M47").Value = 65
$ws.Range("N47").Value = 24000
$ws.Range("O47").Value = 24000
$ws.Range("P47").Value = 24000
$ws.Range("Q47").Value = '$/bandeja 7 kilos'
$ws.Range("S47").Value = 3429
$ws.Range("T47").Value = 7

# Row 48
$ws.Range("D48").Value = 44447
$ws.Range("L48").Value = 'Primera'
$ws.Range("M48").Value = 40
$ws.Range("N48").Value = 3000
$ws.Range("O48").Value = 3000
$ws.Range("P48").Value = 3000
$ws.Range("S48").Value = 3000

# Row 49
$ws.Range("D49").Value = 44187
$ws.Range("L49").Value = 'Primera'
$ws.Range("M49").Value = 15
$ws.Range("N49").Value = 3200
$ws.Range("O49").Value = 3200
$ws.Range("P49").Value = 3200
$ws.Range("S49").Value = 3200

# Row 50
$ws.Range("D50").Value = 44433
$ws.Range("L50").Value = 'Especial'
$ws.Range("N50").Value = 4500
$ws.Range("O50").Value = 4500
$ws.Range("P50").Value = 4500
$ws.Range("S50").Value = 4500

# Row 51
$ws.Range("D51").Value = 44438
$ws.Range("L51").Value = 'Especial'
$ws.Range("M51").Value = 35
$ws.Range("N51").Value = 3500
$ws.Range("O51").Value = 3500
$ws.Range("P51").Value = 3500
$ws.Range("S51").Value = 3500

# Row 52
$ws.Range("D52").Value = 44438
$ws.Range("M52").Value = 20
$ws.Range("N52").Value = 3000
$ws.Range("O52").Value = 3000
$ws.Range("P52").Value = 3000
$ws.Range("S52").Value = 3000

# Row 53
$ws.Range("D53").Value = 44160
$ws.Range("L53").Value = 'Primera'
$ws.Range("M53").Value = 120
$ws.Range("N53").Value = 2200
$ws.Range("O53").Value = 2300
$ws.Range("P53").Value = 2246
$ws.Range("S53").Value = 2246

# Row 54
$ws.Range("D54").Value = 44162
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 85
$ws.Range("N54").Value = 2200
$ws.Range("O54").Value = 2300
$ws.Range("P54").Value = 2247
$ws.Range("S54").Value = 2247

# Row 55
$ws.Range("D55").Value = 44411
$ws.Range("L55").Value = 'Segunda'
$ws.Range("M55").Value = 10
$ws.Range("N55").Value = 3000
$ws.Range("O55").Value = 3000
$ws.Range("P55").Value = 3000
$ws.Range("S55").Value = 3000

# Row 56
$ws.Range("L56").Value = 'Especial'
$ws.Range("M56").Value = 35
$ws.Range("N56").Value = 4500
$ws.Range("O56").Value = 4500
$ws.Range("P56").Value = 4500
$ws.Range("S56").Value = 4500

# --- Append new rows 57-58 (dimension grows from T56 to T58) ---
# Row 57
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = 'Vega Modelo de Temuco'
$ws.Range("C57").Value = 'La Araucanía'
$ws.Range("D57").Value = 44425
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = 'Fruta'
$ws.Range("G57").Value = 100107
$ws.Range("H57").Value = 'Otros'
$ws.Range("I57").Value = 100107002
$ws.Range("J57").Value = 'Chirimoya'
$ws.Range("K57").Value = 'Cultivar IV Región'
$ws.Range("L57").Value = 'Primera'
$ws.Range("M57").Value = 20
$ws.Range("N57").Value = 3500
$ws.Range("O57").Value = 3500
$ws.Range("P57").Value = 3500
$ws.Range("Q57").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R57").Value = 'Provincia del Elquí'
$ws.Range("S57").Value = 3500
$ws.Range("T57").Value = 1

# Row 58
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = 'Vega Modelo de Temuco'
$ws.Range("C58").Value = 'La Araucanía'
$ws.Range("D58").Value = 44425
$ws.Range("D58").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 'Fruta'
$ws.Range("G58").Value = 100107
$ws.Range("H58").Value = 'Otros'
$ws.Range("I58").Value = 100107002
$ws.Range("J58").Value = 'Chirimoya'
$ws.Range("K58").Value = 'Cultivar IV Región'
$ws.Range("L58").Value = 'Segunda'
$ws.Range("M58").Value = 25
$ws.Range("N58").Value = 3000
$ws.Range("O58").Value = 3000
$ws.Range("P58").Value = 3000
$ws.Range("Q58").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R58").Value = 'Provincia del Elquí'
$ws.Range("S58").Value = 3000
$ws.Range("T58").Value = 1
